# Auto-generated edit script: updates computed profit columns (H:N)
# on the leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to
# reflect refreshed market-board pricing data pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 3651.2222  # H6: 6294.2 -> 3651.2222
$ws.Cells.Item(6, 9).Value = 312  # I6: 325.33334 -> 312
$ws.Cells.Item(6, 10).Value = 10329.667  # J6: 15247.5 -> 10329.667
$ws.Cells.Item(6, 11).Value = 936  # K6: 976.0000200000001 -> 936
$ws.Cells.Item(6, 12).Value = 30989.001  # L6: 45742.5 -> 30989.001
$ws.Cells.Item(6, 13).Value = -824  # M6: -864.0000200000001 -> -824
$ws.Cells.Item(6, 14).Value = -31213.001  # N6: -45966.5 -> -31213.001
$ws.Cells.Item(40, 8).Value = 2000  # H40: 2100 -> 2000
$ws.Cells.Item(40, 9).Value = 2000  # I40: 0 -> 2000
$ws.Cells.Item(40, 10).Value = 2000  # J40: 2100 -> 2000
$ws.Cells.Item(40, 11).Value = 2000  # K40: 0 -> 2000
$ws.Cells.Item(40, 12).Value = 2000  # L40: 2100 -> 2000
$ws.Cells.Item(40, 13).Value = -1825  # M40: None -> -1825
$ws.Cells.Item(40, 14).Value = -2350  # N40: -2450 -> -2350
$ws.Cells.Item(51, 8).Value = 5038.095  # H51: 5052.3887 -> 5038.095
$ws.Cells.Item(51, 9).Value = 4953  # I51: 0 -> 4953
$ws.Cells.Item(51, 10).Value = 5042.35  # J51: 5052.3887 -> 5042.35
$ws.Cells.Item(51, 11).Value = 4953  # K51: 0 -> 4953
$ws.Cells.Item(51, 12).Value = 5042.35  # L51: 5052.3887 -> 5042.35
$ws.Cells.Item(51, 13).Value = -4469  # M51: None -> -4469
$ws.Cells.Item(51, 14).Value = -6010.35  # N51: -6020.3887 -> -6010.35
$ws.Cells.Item(127, 8).Value = 3468.7  # H127: 3643.2222 -> 3468.7
$ws.Cells.Item(127, 9).Value = 1632.3334  # I127: 1499.5 -> 1632.3334
$ws.Cells.Item(127, 11).Value = 4897.0002  # K127: 4498.5 -> 4897.0002
$ws.Cells.Item(127, 13).Value = 62.9997999999996  # M127: 461.5 -> 62.9997999999996
$ws.Cells.Item(129, 8).Value = 1786.619  # H129: 11289.81 -> 1786.619
$ws.Cells.Item(129, 9).Value = 924.61536  # I129: 947.5 -> 924.61536
$ws.Cells.Item(129, 10).Value = 3187.375  # J129: 25079.555 -> 3187.375
$ws.Cells.Item(129, 11).Value = 2773.84608  # K129: 2842.5 -> 2773.84608
$ws.Cells.Item(129, 12).Value = 9562.125  # L129: 75238.66500000001 -> 9562.125
$ws.Cells.Item(129, 13).Value = 2226.15392  # M129: 2157.5 -> 2226.15392
$ws.Cells.Item(129, 14).Value = -19562.125  # N129: -85238.66500000001 -> -19562.125
$ws.Cells.Item(131, 8).Value = 2839.6667  # H131: 3029.7273 -> 2839.6667
$ws.Cells.Item(131, 9).Value = 2686.7778  # I131: 2929 -> 2686.7778
$ws.Cells.Item(131, 11).Value = 8060.3334  # K131: 8787 -> 8060.3334
$ws.Cells.Item(131, 13).Value = -3020.3334  # M131: -3747 -> -3020.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5930.9316  # H32: 6059.6514 -> 5930.9316
$ws.Cells.Item(32, 9).Value = 5318.5  # I32: 5459.143 -> 5318.5
$ws.Cells.Item(32, 11).Value = 5318.5  # K32: 5459.143 -> 5318.5
$ws.Cells.Item(32, 13).Value = -5031.5  # M32: -5172.143 -> -5031.5
$ws.Cells.Item(37, 8).Value = 39940  # H37: 0 -> 39940
$ws.Cells.Item(37, 10).Value = 39940  # J37: 0 -> 39940
$ws.Cells.Item(37, 12).Value = 39940  # L37: 0 -> 39940
$ws.Cells.Item(37, 14).Value = -40486  # N37: None -> -40486
$ws.Cells.Item(45, 8).Value = 4356  # H45: 5049.3335 -> 4356
$ws.Cells.Item(45, 9).Value = 5478  # I45: 6249.4287 -> 5478
$ws.Cells.Item(45, 10).Value = 3374.25  # J45: 3999.25 -> 3374.25
$ws.Cells.Item(45, 11).Value = 5478  # K45: 6249.4287 -> 5478
$ws.Cells.Item(45, 12).Value = 3374.25  # L45: 3999.25 -> 3374.25
$ws.Cells.Item(45, 13).Value = -5101  # M45: -5872.4287 -> -5101
$ws.Cells.Item(45, 14).Value = -4128.25  # N45: -4753.25 -> -4128.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(39, 8).Value = 15000  # H39: 0 -> 15000
$ws.Cells.Item(39, 10).Value = 15000  # J39: 0 -> 15000
$ws.Cells.Item(39, 12).Value = 15000  # L39: 0 -> 15000
$ws.Cells.Item(39, 14).Value = -15778  # N39: None -> -15778
$ws.Cells.Item(86, 8).Value = 49657.1  # H86: 61252.625 -> 49657.1
$ws.Cells.Item(86, 9).Value = 36446.375  # I86: 47503.5 -> 36446.375
$ws.Cells.Item(86, 11).Value = 36446.375  # K86: 47503.5 -> 36446.375
$ws.Cells.Item(86, 13).Value = -35323.375  # M86: -46380.5 -> -35323.375
$ws.Cells.Item(89, 8).Value = 49657.1  # H89: 61252.625 -> 49657.1
$ws.Cells.Item(89, 9).Value = 36446.375  # I89: 47503.5 -> 36446.375
$ws.Cells.Item(89, 11).Value = 182231.875  # K89: 237517.5 -> 182231.875
$ws.Cells.Item(89, 13).Value = -176615.875  # M89: -231901.5 -> -176615.875
$ws.Cells.Item(94, 8).Value = 1616.8096  # H94: 1619.238 -> 1616.8096
$ws.Cells.Item(94, 9).Value = 1462.1111  # I94: 1467.7778 -> 1462.1111
$ws.Cells.Item(94, 11).Value = 1462.1111  # K94: 1467.7778 -> 1462.1111
$ws.Cells.Item(94, 13).Value = -1011.1111  # M94: -1016.7778 -> -1011.1111

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1015.4286  # H16: 1284.8334 -> 1015.4286
$ws.Cells.Item(16, 9).Value = 1015.4286  # I16: 1284.8334 -> 1015.4286
$ws.Cells.Item(16, 11).Value = 1015.4286  # K16: 1284.8334 -> 1015.4286
$ws.Cells.Item(16, 13).Value = -728.4286  # M16: -997.8334 -> -728.4286
$ws.Cells.Item(31, 8).Value = 11909751  # H31: 14711151 -> 11909751
$ws.Cells.Item(31, 9).Value = 2831.5334  # I31: 3044.2307 -> 2831.5334
$ws.Cells.Item(31, 10).Value = 41677050  # J31: 62512500 -> 41677050
$ws.Cells.Item(31, 11).Value = 2831.5334  # K31: 3044.2307 -> 2831.5334
$ws.Cells.Item(31, 12).Value = 41677050  # L31: 62512500 -> 41677050
$ws.Cells.Item(31, 13).Value = -2536.5334  # M31: -2749.2307 -> -2536.5334
$ws.Cells.Item(31, 14).Value = -41677640  # N31: -62513090 -> -41677640
$ws.Cells.Item(34, 8).Value = 11909751  # H34: 14711151 -> 11909751
$ws.Cells.Item(34, 9).Value = 2831.5334  # I34: 3044.2307 -> 2831.5334
$ws.Cells.Item(34, 10).Value = 41677050  # J34: 62512500 -> 41677050
$ws.Cells.Item(34, 11).Value = 2831.5334  # K34: 3044.2307 -> 2831.5334
$ws.Cells.Item(34, 12).Value = 41677050  # L34: 62512500 -> 41677050
$ws.Cells.Item(34, 13).Value = -2629.5334  # M34: -2842.2307 -> -2629.5334
$ws.Cells.Item(34, 14).Value = -41677454  # N34: -62512904 -> -41677454
$ws.Cells.Item(54, 8).Value = 39724.25  # H54: 39724.75 -> 39724.25
$ws.Cells.Item(54, 10).Value = 37998.5  # J54: 37999.5 -> 37998.5
$ws.Cells.Item(54, 12).Value = 37998.5  # L54: 37999.5 -> 37998.5
$ws.Cells.Item(54, 14).Value = -39314.5  # N54: -39315.5 -> -39314.5
$ws.Cells.Item(86, 8).Value = 0  # H86: 3999.5 -> 0
$ws.Cells.Item(86, 9).Value = 0  # I86: 3999.5 -> 0
$ws.Cells.Item(86, 11).Value = 0  # K86: 3999.5 -> 0
$ws.Cells.Item(86, 13).ClearContents()  # M86: -2876.5 -> (blank)
$ws.Cells.Item(89, 8).Value = 0  # H89: 3999.5 -> 0
$ws.Cells.Item(89, 9).Value = 0  # I89: 3999.5 -> 0
$ws.Cells.Item(89, 11).Value = 0  # K89: 19997.5 -> 0
$ws.Cells.Item(89, 13).ClearContents()  # M89: -14381.5 -> (blank)
$ws.Cells.Item(113, 8).Value = 1015.4286  # H113: 1284.8334 -> 1015.4286
$ws.Cells.Item(113, 9).Value = 1015.4286  # I113: 1284.8334 -> 1015.4286
$ws.Cells.Item(113, 11).Value = 1015.4286  # K113: 1284.8334 -> 1015.4286
$ws.Cells.Item(113, 13).Value = 1154.5714  # M113: 885.1666 -> 1154.5714
$ws.Cells.Item(141, 8).Value = 141598.8  # H141: 121080.57 -> 141598.8
$ws.Cells.Item(141, 10).Value = 169748.5  # J141: 136427.33 -> 169748.5
$ws.Cells.Item(141, 12).Value = 169748.5  # L141: 136427.33 -> 169748.5
$ws.Cells.Item(141, 14).Value = -180108.5  # N141: -146787.33 -> -180108.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 22378900  # H4: 25313826 -> 22378900
$ws.Cells.Item(4, 9).Value = 23075900  # I4: 27809398 -> 23075900
$ws.Cells.Item(4, 11).Value = 69227700  # K4: 83428194 -> 69227700
$ws.Cells.Item(4, 13).Value = -69227588  # M4: -83428082 -> -69227588
$ws.Cells.Item(23, 8).Value = 865  # H23: 957.6 -> 865
$ws.Cells.Item(23, 9).Value = 273.8  # I23: 459 -> 273.8
$ws.Cells.Item(23, 11).Value = 821.4000000000001  # K23: 1377 -> 821.4000000000001
$ws.Cells.Item(23, 13).Value = -586.4000000000001  # M23: -1142 -> -586.4000000000001
$ws.Cells.Item(34, 8).Value = 1692.0385  # H34: 1662.7037 -> 1692.0385
$ws.Cells.Item(34, 10).Value = 4746.1113  # J34: 4361.5 -> 4746.1113
$ws.Cells.Item(34, 12).Value = 14238.3339  # L34: 13084.5 -> 14238.3339
$ws.Cells.Item(34, 14).Value = -14406.3339  # N34: -13252.5 -> -14406.3339
$ws.Cells.Item(39, 8).Value = 673.9286  # H39: 902.125 -> 673.9286
$ws.Cells.Item(39, 9).Value = 455.3846  # I39: 529.9286 -> 455.3846
$ws.Cells.Item(39, 10).Value = 3515  # J39: 3507.5 -> 3515
$ws.Cells.Item(39, 11).Value = 1366.1538  # K39: 1589.7858 -> 1366.1538
$ws.Cells.Item(39, 12).Value = 10545  # L39: 10522.5 -> 10545
$ws.Cells.Item(39, 13).Value = -1072.1538  # M39: -1295.7858 -> -1072.1538
$ws.Cells.Item(39, 14).Value = -11133  # N39: -11110.5 -> -11133
$ws.Cells.Item(55, 8).Value = 1005177.25  # H55: 905077 -> 1005177.25
$ws.Cells.Item(55, 10).Value = 6099.2856  # J55: 5858.75 -> 6099.2856
$ws.Cells.Item(55, 12).Value = 18297.8568  # L55: 17576.25 -> 18297.8568
$ws.Cells.Item(55, 14).Value = -18651.8568  # N55: -17930.25 -> -18651.8568
$ws.Cells.Item(131, 8).Value = 20001164  # H131: 20001174 -> 20001164
$ws.Cells.Item(131, 9).Value = 33334108  # I131: 33334124 -> 33334108
$ws.Cells.Item(131, 11).Value = 100002324  # K131: 100002372 -> 100002324
$ws.Cells.Item(131, 13).Value = -99997284  # M131: -99997332 -> -99997284
$ws.Cells.Item(134, 8).Value = 1054.6154  # H134: 4799.2666 -> 1054.6154
$ws.Cells.Item(134, 9).Value = 889.75  # I134: 814.1818 -> 889.75
$ws.Cells.Item(134, 10).Value = 3033  # J134: 15758.25 -> 3033
$ws.Cells.Item(134, 11).Value = 2669.25  # K134: 2442.5454 -> 2669.25
$ws.Cells.Item(134, 12).Value = 9099  # L134: 47274.75 -> 9099
$ws.Cells.Item(134, 13).Value = 2400.75  # M134: 2627.4546 -> 2400.75
$ws.Cells.Item(134, 14).Value = -19239  # N134: -57414.75 -> -19239

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 4428.7144  # H80: 4658.231 -> 4428.7144
$ws.Cells.Item(80, 10).Value = 6211.5  # J80: 6892.4287 -> 6211.5
$ws.Cells.Item(80, 12).Value = 6211.5  # L80: 6892.4287 -> 6211.5
$ws.Cells.Item(80, 14).Value = -8207.5  # N80: -8888.4287 -> -8207.5
$ws.Cells.Item(83, 8).Value = 4428.7144  # H83: 4658.231 -> 4428.7144
$ws.Cells.Item(83, 10).Value = 6211.5  # J83: 6892.4287 -> 6211.5
$ws.Cells.Item(83, 12).Value = 31057.5  # L83: 34462.14350000001 -> 31057.5
$ws.Cells.Item(83, 14).Value = -41041.5  # N83: -44446.14350000001 -> -41041.5
$ws.Cells.Item(97, 8).Value = 1050.9584  # H97: 945.5185 -> 1050.9584
$ws.Cells.Item(97, 9).Value = 1173.125  # I97: 993.4 -> 1173.125
$ws.Cells.Item(97, 10).Value = 806.625  # J97: 808.7143 -> 806.625
$ws.Cells.Item(97, 11).Value = 1173.125  # K97: 993.4 -> 1173.125
$ws.Cells.Item(97, 12).Value = 806.625  # L97: 808.7143 -> 806.625
$ws.Cells.Item(97, 13).Value = -677.125  # M97: -497.4 -> -677.125
$ws.Cells.Item(97, 14).Value = -1798.625  # N97: -1800.7143 -> -1798.625
$ws.Cells.Item(132, 8).Value = 2060.7334  # H132: 2023.2903 -> 2060.7334
$ws.Cells.Item(132, 9).Value = 2053.48  # I132: 2009.1154 -> 2053.48
$ws.Cells.Item(132, 11).Value = 6160.440000000001  # K132: 6027.3462 -> 6160.440000000001
$ws.Cells.Item(132, 13).Value = -3630.440000000001  # M132: -3497.3462 -> -3630.440000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2842.2307  # H22: 2960.6428 -> 2842.2307
$ws.Cells.Item(22, 9).Value = 2194.9  # I22: 2444.9 -> 2194.9
$ws.Cells.Item(22, 10).Value = 5000  # J22: 4250 -> 5000
$ws.Cells.Item(22, 11).Value = 2194.9  # K22: 2444.9 -> 2194.9
$ws.Cells.Item(22, 12).Value = 5000  # L22: 4250 -> 5000
$ws.Cells.Item(22, 13).Value = -1899.9  # M22: -2149.9 -> -1899.9
$ws.Cells.Item(22, 14).Value = -5590  # N22: -4840 -> -5590
$ws.Cells.Item(27, 8).Value = 2842.2307  # H27: 2960.6428 -> 2842.2307
$ws.Cells.Item(27, 9).Value = 2194.9  # I27: 2444.9 -> 2194.9
$ws.Cells.Item(27, 10).Value = 5000  # J27: 4250 -> 5000
$ws.Cells.Item(27, 11).Value = 2194.9  # K27: 2444.9 -> 2194.9
$ws.Cells.Item(27, 12).Value = 5000  # L27: 4250 -> 5000
$ws.Cells.Item(27, 13).Value = -2087.9  # M27: -2337.9 -> -2087.9
$ws.Cells.Item(27, 14).Value = -5214  # N27: -4464 -> -5214
$ws.Cells.Item(55, 8).Value = 908  # H55: 932.7619 -> 908
$ws.Cells.Item(55, 9).Value = 537  # I55: 602.7 -> 537
$ws.Cells.Item(55, 10).Value = 1316.1  # J55: 1232.8182 -> 1316.1
$ws.Cells.Item(55, 11).Value = 537  # K55: 602.7 -> 537
$ws.Cells.Item(55, 12).Value = 1316.1  # L55: 1232.8182 -> 1316.1
$ws.Cells.Item(55, 13).Value = -364  # M55: -429.7 -> -364
$ws.Cells.Item(55, 14).Value = -1662.1  # N55: -1578.8182 -> -1662.1
$ws.Cells.Item(82, 8).Value = 2149.8333  # H82: 2022.3846 -> 2149.8333
$ws.Cells.Item(82, 9).Value = 2459.2  # I82: 2116 -> 2459.2
$ws.Cells.Item(82, 10).Value = 1928.8572  # J82: 1942.1428 -> 1928.8572
$ws.Cells.Item(82, 11).Value = 2459.2  # K82: 2116 -> 2459.2
$ws.Cells.Item(82, 12).Value = 1928.8572  # L82: 1942.1428 -> 1928.8572
$ws.Cells.Item(82, 13).Value = -2098.2  # M82: -1755 -> -2098.2
$ws.Cells.Item(82, 14).Value = -2650.8572  # N82: -2664.1428 -> -2650.8572
$ws.Cells.Item(85, 8).Value = 2149.8333  # H85: 2022.3846 -> 2149.8333
$ws.Cells.Item(85, 9).Value = 2459.2  # I85: 2116 -> 2459.2
$ws.Cells.Item(85, 10).Value = 1928.8572  # J85: 1942.1428 -> 1928.8572
$ws.Cells.Item(85, 11).Value = 2459.2  # K85: 2116 -> 2459.2
$ws.Cells.Item(85, 12).Value = 1928.8572  # L85: 1942.1428 -> 1928.8572
$ws.Cells.Item(85, 13).Value = -1211.2  # M85: -868 -> -1211.2
$ws.Cells.Item(85, 14).Value = -4424.8572  # N85: -4438.1428 -> -4424.8572
$ws.Cells.Item(131, 8).Value = 89977  # H131: 89976.25 -> 89977
$ws.Cells.Item(131, 10).Value = 89977  # J131: 89976.25 -> 89977
$ws.Cells.Item(131, 12).Value = 89977  # L131: 89976.25 -> 89977
$ws.Cells.Item(131, 14).Value = -100057  # N131: -100056.25 -> -100057
$ws.Cells.Item(132, 8).Value = 6762.304  # H132: 7225.476 -> 6762.304
$ws.Cells.Item(132, 9).Value = 2929.1052  # I132: 3050.2942 -> 2929.1052
$ws.Cells.Item(132, 11).Value = 8787.3156  # K132: 9150.882599999999 -> 8787.3156
$ws.Cells.Item(132, 13).Value = -6257.3156  # M132: -6620.882599999999 -> -6257.3156

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(8, 8).Value = 3666.5  # H8: 3666 -> 3666.5
$ws.Cells.Item(8, 9).Value = 3666.5  # I8: 3666 -> 3666.5
$ws.Cells.Item(8, 11).Value = 3666.5  # K8: 3666 -> 3666.5
$ws.Cells.Item(8, 13).Value = -3526.5  # M8: -3526 -> -3526.5
$ws.Cells.Item(81, 8).Value = 2435.0908  # H81: 2460.4546 -> 2435.0908
$ws.Cells.Item(81, 9).Value = 1112.4286  # I81: 1152.2858 -> 1112.4286
$ws.Cells.Item(81, 11).Value = 2224.8572  # K81: 2304.5716 -> 2224.8572
$ws.Cells.Item(81, 13).Value = -1163.8572  # M81: -1243.5716 -> -1163.8572
$ws.Cells.Item(84, 8).Value = 2435.0908  # H84: 2460.4546 -> 2435.0908
$ws.Cells.Item(84, 9).Value = 1112.4286  # I84: 1152.2858 -> 1112.4286
$ws.Cells.Item(84, 11).Value = 11124.286  # K84: 11522.858 -> 11124.286
$ws.Cells.Item(84, 13).Value = -5820.286  # M84: -6218.858 -> -5820.286
$ws.Cells.Item(132, 8).Value = 2011.0244  # H132: 2066.878 -> 2011.0244
$ws.Cells.Item(132, 9).Value = 1493.8276  # I132: 1572.7931 -> 1493.8276
$ws.Cells.Item(132, 11).Value = 4481.4828  # K132: 4718.379300000001 -> 4481.4828
$ws.Cells.Item(132, 13).Value = -1951.4828  # M132: -2188.379300000001 -> -1951.4828
$ws.Cells.Item(136, 8).Value = 9659.482  # H136: 9897.321 -> 9659.482
$ws.Cells.Item(136, 9).Value = 3514.5  # I136: 3551.8462 -> 3514.5
$ws.Cells.Item(136, 10).Value = 13997.117  # J136: 15396.733 -> 13997.117
$ws.Cells.Item(136, 11).Value = 10543.5  # K136: 10655.5386 -> 10543.5
$ws.Cells.Item(136, 12).Value = 41991.351  # L136: 46190.199 -> 41991.351
$ws.Cells.Item(136, 13).Value = -7993.5  # M136: -8105.5386 -> -7993.5
$ws.Cells.Item(136, 14).Value = -47091.351  # N136: -51290.199 -> -47091.351
